$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# families_descriptions sheet: insert a new "title" column (B) between
# the family code (A) and the description (shifts from B to C).
# ------------------------------------------------------------------
$fd = $wb.Worksheets.Item("families_descriptions")

$codeColumnWidth = $fd.Columns("A:A").ColumnWidth

$fd.Columns("B:B").Insert()
$fd.Columns("B:B").ColumnWidth = $codeColumnWidth

$fd.Range("A1").Value = "code"
$fd.Range("B1").Value = "label-en_US"
$fd.Range("C1").Value = "description-en_US"

$fd.Range("B2").Value = "Videogames"
$fd.Range("B3").Value = "Jerseys"
$fd.Range("B4").Value = "Headphones"

# ------------------------------------------------------------------
# Workbook / sheet view state: families_descriptions becomes the
# active tab (instead of attribute_options), with cell C8 selected.
# ------------------------------------------------------------------
$fd.Activate()
$fd.Range("C8").Select()
